$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after testdata_Prop and rename it
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "testdata_Prop_g"

# Copy header row (A1:J1) from testdata_Prop, keeping formatting (styles s=1 / s=3)
$ws1.Range("A1:J1").Copy($newSheet.Range("A1:J1"))

# Copy the Area-name column (A) and the constant metadata columns (G:J) with formatting
$ws1.Range("A2:A9").Copy($newSheet.Range("A2:A9"))
$ws1.Range("G2:J9").Copy($newSheet.Range("G2:J9"))

# Column B/C: SUMIF totals of the matching Area rows on testdata_Prop
for ($r = 2; $r -le 9; $r++) {
    $newSheet.Cells.Item($r, 2).Formula = "=SUMIF(testdata_Prop!`$A`$2:`$A`$33,testdata_Prop_g!`$A$r,testdata_Prop!B`$2:B`$33)"
    $newSheet.Cells.Item($r, 3).Formula = "=SUMIF(testdata_Prop!`$A`$2:`$A`$33,testdata_Prop_g!`$A$r,testdata_Prop!C`$2:C`$33)"
}

# Column D/E/F: proportion / lowercl / uppercl values
# (written in plain-decimal form -- the PowerShell parser here doesn't
# accept scientific-notation numeric literals like 1.23E-4)
$newSheet.Range("D2").Value = 0.79837991089509919
$newSheet.Range("E2").Value = 0.79481805800067218
$newSheet.Range("F2").Value = 0.80189534317473543

$newSheet.Range("D3").Value = 0.05
$newSheet.Range("E3").Value = 0.032597429837147258
$newSheet.Range("F3").Value = 0.075963635063719587

$newSheet.Range("D4").Value = 0.2
$newSheet.Range("E4").Value = 0.16373705973387687
$newSheet.Range("F4").Value = 0.24197031686670104

$newSheet.Range("D5").Value = 0.65
$newSheet.Range("E5").Value = 0.60203196070747378
$newSheet.Range("F5").Value = 0.69511435099223728

$newSheet.Range("D6").Value = 0.79837991089509919
$newSheet.Range("E6").Value = 0.79481805800067218
$newSheet.Range("F6").Value = 0.80189534317473543

$newSheet.Range("D7").Value = 0.0010280420256339154
$newSheet.Range("E7").Value = 0.0010276809995663644
$newSheet.Range("F7").Value = 0.0010284031783999669

$newSheet.Range("D8").Value = 0.069051321928460335
$newSheet.Range("E8").Value = 0.064796737148613115
$newSheet.Range("F8").Value = 0.073563290415335325

$newSheet.Range("D9").Value = 0.2890625
$newSheet.Range("E9").Value = 0.28156503717401421
$newSheet.Range("F9").Value = 0.29667716227788915

# Column C is a little wider on the new sheet (to fit the big denominators)
$newSheet.Columns.Item(3).ColumnWidth = 11.6

# Selections: testdata_Prop keeps a (now inactive) selection of G2:J7,
# testdata_Prop_g becomes the active tab with E26 selected
$ws1.Range("G2:J7").Select()
$newSheet.Range("E26").Select()
